# Updated cryptos list values (prices + 1h volume change %) per upstream diff.
# Source sheet stores these as plain text (t="inlineStr"/"s"), so numeric-looking
# price strings are forced to Text before assignment to avoid Excel silently
# coercing them to numbers (which would drop formatting like trailing zeros,
# e.g. "0.5200" -> 0.52) and restored to the default "Normal" style afterwards so
# no stray per-cell style/quote-prefix marker is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.498.40"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "1.845.46"
$ws.Range("E3").Value = "  -0.19%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "263.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5200"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3212"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06779"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7712"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07774"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "1.859.64"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.014"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("E17").Value = "  -1.48%  "
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007940"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.58%  "
$ws.Range("D20").Value = "26.530.83"
$ws.Range("E20").Value = "  +0.06%  "
$ws.Range("D21").Value = "2.089.46"
$ws.Range("E21").Value = "  -0.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.610"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.429"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.981"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.28%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.177"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.84%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.677"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.97"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "111.38"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.159"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08724"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.098"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04808"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.68%  "
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7183"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.91%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.862"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.094"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01782"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.196"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.78%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.4827"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "112.01"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8951"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.032"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.000"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.613"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4164"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05901"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.012"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1227"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.25%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8850"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.74%  "
